# Weekly refresh of the Ciboulette price series:
# two brand-new rows (most recent week) are inserted at the top of the
# data block and every older row shifts down by two, with the oldest
# two rows falling off the bottom of the previously-used range (so the
# sheet grows from 44 to 46 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 38..46 (row 37 is untouched).
# Columns A,B,C,E,F,G,H,N,O,Q,R are constant across the whole sheet.
$constA = 7
$constB = "Terminal Hortofrutícola Agro Chillán"
$constC = "Ñuble"
$constE = 16
$constF = 100112039
$constG = "Ciboulette"
$constH = "Sin especificar"
$constN = "`$/docena de atados"
$constO = "Región Metropolitana"
$constQ = 3
$constR = "Hortaliza"

$rows = @(
    @{ Row=38; D=45258; I="Primera"; J=100; K=2500; L=2500; M=2500; P=833 },
    @{ Row=39; D=45258; I="Segunda"; J=100; K=2000; L=2000; M=2000; P=667 },
    @{ Row=40; D=45233; I="Primera"; J=300; K=2000; L=2000; M=2000; P=667 },
    @{ Row=41; D=45233; I="Segunda"; J=250; K=1500; L=1500; M=1500; P=500 },
    @{ Row=42; D=45145; I="Primera"; J=60;  K=2500; L=2500; M=2500; P=833 },
    @{ Row=43; D=45145; I="Segunda"; J=80;  K=2000; L=2000; M=2000; P=667 },
    @{ Row=44; D=44832; I="Primera"; J=200; K=1200; L=1300; M=1250; P=417 },
    @{ Row=45; D=44832; I="Segunda"; J=150; K=1000; L=1000; M=1000; P=333 },
    @{ Row=46; D=45135; I="Primera"; J=70;  K=2500; L=2500; M=2500; P=833 }
)

$dateFormat = $ws.Range("D37").NumberFormat

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $constA
    $ws.Cells.Item($row, 2).Value = $constB
    $ws.Cells.Item($row, 3).Value = $constC

    $ws.Cells.Item($row, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 4).Value = $r.D

    $ws.Cells.Item($row, 5).Value = $constE
    $ws.Cells.Item($row, 6).Value = $constF
    $ws.Cells.Item($row, 7).Value = $constG
    $ws.Cells.Item($row, 8).Value = $constH
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $constN
    $ws.Cells.Item($row, 15).Value = $constO
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $constQ
    $ws.Cells.Item($row, 18).Value = $constR
}

Write-Host "Updated rows 38-46"
